$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column G ("k") - shifts old G:O (ARI..Predicted_Subgraphs) to H:P
$ws.Columns("G:G").Insert(-4161)

# Insert two new columns J:K ("GNMI","FuzzyARI") - shifts old I:P (Purity..Predicted_Subgraphs) to L:R
$ws.Columns("J:K").Insert(-4161)

# --- Header row (row 1) ---
# Columns inserted via Insert() already inherit the bold/border/center header
# style from their neighboring header cells, so plain value assignment is enough.
$ws.Range("G1").Value = "k"
$ws.Range("J1").Value = "GNMI"
$ws.Range("K1").Value = "FuzzyARI"

# --- Row 2 (Graph1) ---
$ws.Range("G2").Value = 0.001
$ws.Range("H2").Value = 0.09469922572960095
$ws.Range("I2").Value = 0.3107373603223778
$ws.Range("J2").Value = 0.3099513282337183
$ws.Range("K2").Value = 0.09469922572960096
$ws.Range("L2").Value = 0.45
$ws.Range("M2").Value = 2.270950594454669
$ws.Range("N2").Value = 0.6098403047164005
$ws.Range("O2").Value = 0.2314285714285714
$ws.Range("P2").Value = 0.4429522220587847
$ws.Range("Q2").Value = "Subgraph 1 (Nodes): [20, 25, 21] - Density: 1`nSubgraph 2 (Nodes): [1, 4, 24, 11, 29, 15] - Density: 0.799943`nSubgraph 3 (Nodes): [10, 7, 2] - Density: 1`nSubgraph 4 (Nodes): [14, 8, 0, 19] - Density: 1`n----------------------------------------------------"
$ws.Range("R2").Value = "Subgraph 1:z { 0 1 2 4 6 7 8 10 11 14 15 19 20 21 24 25 29 } N: 17 Triangles: 23 Density: 0.0338178"

# --- Row 3 (Graph2) ---
$ws.Range("G3").Value = 0.001
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0.2340425531914894
$ws.Range("M3").Value = 2.526255260617095
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0.05357142857142858
$ws.Range("P3").Value = 0.4023520119691402
$ws.Range("Q3").Value = "Subgraph 1 (Nodes): [27, 20, 17, 37, 28, 30, 12, 49, 39] - Density: 0.761897`nSubgraph 2 (Nodes): [45, 34, 21, 41, 26] - Density: 1`nSubgraph 3 (Nodes): [0, 7, 40, 43, 44] - Density: 1`nSubgraph 4 (Nodes): [33, 22, 2, 32, 42, 18, 16, 38] - Density: 0.80358`nSubgraph 5 (Nodes): [6, 15, 13, 11, 48, 3, 10, 47, 29] - Density: 0.833321`n----------------------------------------------------"
$ws.Range("R3").Value = "Subgraph 1:z { 0 1 2 3 5 6 7 8 9 10 11 12 13 14 15 16 17 18 19 20 21 22 24 25 26 27 28 29 30 31 32 33 34 35 37 38 39 40 41 42 43 44 45 46 47 48 49 } N: 47 Triangles: 199 Density: 0.0122144"

# Setting multi-line text triggers Excel's wrap-text auto row height; reset it
# back to the sheet default so row heights match the un-customized original.
$ws.Rows("2:3").EntireRow.AutoFit()
